$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from E1 (existing bold/border style) onto F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Header label for the new column
$ws.Range("F1").Value = "time_taken"

# Per-row "time_taken" metadata timestamps (microsecond precision strings)
$timestamps = @(
    "2021-10-05 13:39:26.814090",
    "2021-10-05 13:39:26.814104",
    "2021-10-05 13:39:26.814109",
    "2021-10-05 13:39:26.814112",
    "2021-10-05 13:39:26.814115",
    "2021-10-05 13:39:26.814119",
    "2021-10-05 13:39:26.814122",
    "2021-10-05 13:39:26.814125",
    "2021-10-05 13:39:26.814128",
    "2021-10-05 13:39:26.814131",
    "2021-10-05 13:39:26.814134",
    "2021-10-05 13:39:26.814137",
    "2021-10-05 13:39:26.814140",
    "2021-10-05 13:39:26.814143",
    "2021-10-05 13:39:26.814146",
    "2021-10-05 13:39:26.814149",
    "2021-10-05 13:39:26.814152",
    "2021-10-05 13:39:26.814156",
    "2021-10-05 13:39:26.814159",
    "2021-10-05 13:39:26.814162",
    "2021-10-05 13:39:26.814165",
    "2021-10-05 13:39:26.814168",
    "2021-10-05 13:39:26.814171",
    "2021-10-05 13:39:26.814174",
    "2021-10-05 13:39:26.814178",
    "2021-10-05 13:39:26.814181",
    "2021-10-05 13:39:26.814184",
    "2021-10-05 13:39:26.814187",
    "2021-10-05 13:39:26.814190",
    "2021-10-05 13:39:26.814193",
    "2021-10-05 13:39:26.814196",
    "2021-10-05 13:39:26.814199",
    "2021-10-05 13:39:26.814203",
    "2021-10-05 13:39:26.814206"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Host "Added time_taken column F1:F35"
